# Generate Report for Handoff
# Update the "Latest Handoff Datetime" (column D, row 5) on both the
# zh-cn and de-de report sheets to reflect the new handoff timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Cells.Item(5, 4).Value = "2016-02-22 05:55:04"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Cells.Item(5, 4).Value = "2016-02-22 05:55:17"
